# Insert a new data row above row 389 (pushing existing rows 389..482 down to 390..483)
# and populate the newly inserted row with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(389).Insert()

$ws.Range("A389").Value = 9
$ws.Range("B389").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C389").Value = "Metropolitana"
$ws.Range("D389").Value = 44798
$ws.Range("E389").Value = 13
$ws.Range("F389").Value = 100112012
$ws.Range("G389").Value = "Espinaca"
$ws.Range("H389").Value = "Sin especificar"
$ws.Range("I389").Value = "Primera"
$ws.Range("J389").Value = 160
$ws.Range("K389").Value = 6000
$ws.Range("L389").Value = 7000
$ws.Range("M389").Value = 6500
$ws.Range("N389").Value = "`$/cuna 10 kilos"
$ws.Range("O389").Value = "Provincia de Chacabuco"
$ws.Range("P389").Value = 650
$ws.Range("Q389").Value = 10
$ws.Range("R389").Value = "Hortaliza"
